$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$malfegor = "('Malfegor', ['{2}{B}{B}{R}{R}', 'Legendary Creature " + [char]0x2014 + " Demon Dragon', 'Flying', 'When Malfegor enters the battlefield, discard your hand. Each opponent sacrifices a creature for each card discarded this way.', '6/6'])"
$obelisk = "('Obelisk of Alara', ['{6}', 'Artifact', '{1}{W}, {T}: You gain 5 life.', '{1}{U}, {T}: Draw a card, then discard a card.', '{1}{B}, {T}: Target creature gets -2/-2 until end of turn.', '{1}{R}, {T}: Obelisk of Alara deals 3 damage to target player or planeswalker.', '{1}{G}, {T}: Target creature gets +4/+4 until end of turn.'])"

$ws.Range("A2").Value = $malfegor
$ws.Range("A3").Value = $obelisk

# Remove old rows 4 through 15, which are no longer needed
$ws.Range("A4:A15").EntireRow.Delete()
